$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "BrowserUtil" (sheet1.xml)
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("BrowserUtil")

# Remove existing hyperlinks/extra row before rebuilding the table
$ws1.Hyperlinks.Delete()
$ws1.Rows.Item(9).Delete()

$ws1Data = @(
    @("TestName",               "URL",   "UserName"),
    @("bc_SelectAnEnvironment",  "URL2",  "Loadt02"),
    @("bc_Login",                "URL3",  "Loadt03"),
    @("bc_ClickOnEnter",         "URL4",  "Loadt04"),
    @("bc_SearchCode",           "URL5",  "Loadt05"),
    @("tc_LoginTest",            "https://aenetworks.oktapreview.com/login/default", "Loadt06"),
    @("bc_FIDocsUpload",         "URL7",  "Loadt07"),
    @("bc_ClickOnExecute",       "URL8",  "Loadt08")
)

for ($i = 0; $i -lt $ws1Data.Length; $i++) {
    $r = $i + 1
    $row = $ws1Data[$i]
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
}

# Hyperlinks (added in order so relationship ids come out rId1..rId4).
# Adding a hyperlink with a "TextToDisplay" overwrites the cell's text with
# that display text, so the original placeholder text is restored right after.
$ws1.Hyperlinks.Add($ws1.Range("B3"), "https://aenetworks.oktapreview.com/login/default", "", "", "https://aenetworks.oktapreview.com/login/default") | Out-Null
$ws1.Range("B3").Value = "URL3"
$ws1.Hyperlinks.Add($ws1.Range("B5"), "https://aenetworks.oktapreview.com/login/default", "", "", "https://aenetworks.oktapreview.com/login/default") | Out-Null
$ws1.Range("B5").Value = "URL5"
$ws1.Hyperlinks.Add($ws1.Range("B7"), "https://aenetworks.oktapreview.com/login/default", "", "", "https://aenetworks.oktapreview.com/login/default") | Out-Null
$ws1.Range("B7").Value = "URL7"
$ws1.Hyperlinks.Add($ws1.Range("B6"), "https://aenetworks.oktapreview.com/login/default") | Out-Null

# Styles: column A (rows 2-8) yellow fill, column B (rows 2-8) hyperlink look.
# Re-applied after the hyperlinks are created so every cell ends up sharing
# the same style record instead of Excel creating per-cell variants.
$ws1.Range("A2:A8").Interior.Color = 65535
$ws1.Range("B2:B8").Style = "Hyperlink"

$ws1.Activate() | Out-Null
$ws1.Range("A5").Select() | Out-Null

# ---------------------------------------------------------------
# Sheet "Sheet1" (sheet2.xml)
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet1")

$ws2Data = @(
    @("TestName",               "URL",  "UserName"),
    @("tc_LoginTest",           "https://aenetworks.oktapreview.com/login/default", "Loadt01"),
    @("bc_SelectAnEnvironment", "https://aenetworks.oktapreview.com/login/default", "Loadt02"),
    @("bc_Login",               "https://aenetworks.oktapreview.com/login/default", "Loadt03"),
    @("bc_ClickOnEnter",        "https://aenetworks.oktapreview.com/login/default", "Loadt04"),
    @("bc_SearchCode",          "https://aenetworks.oktapreview.com/login/default", "Loadt05"),
    @("bc_ClickOnEnter",        "https://aenetworks.oktapreview.com/login/default", "Loadt06"),
    @("bc_FIDocsUpload",        "https://aenetworks.oktapreview.com/login/default", "Loadt07"),
    @("bc_ClickOnExecute",      "https://aenetworks.oktapreview.com/login/default", "Loadt08")
)

for ($i = 0; $i -lt $ws2Data.Length; $i++) {
    $r = $i + 1
    $row = $ws2Data[$i]
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
}

$ws2.Columns.Item(1).ColumnWidth = 22.833333333333336
$ws2.Columns.Item(2).ColumnWidth = 47.33333333333333

$ws2.Hyperlinks.Add($ws2.Range("B2"), "https://aenetworks.oktapreview.com/login/default") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B3:B9"), "https://aenetworks.oktapreview.com/login/default", "", "", "https://aenetworks.oktapreview.com/login/default") | Out-Null

# Re-apply styling after hyperlinks are created so every cell shares the
# same style record instead of Excel creating per-cell style variants.
$ws2.Range("A2:A9").Interior.Color = 65535
$ws2.Range("B2:B9").Style = "Hyperlink"

$ws2.Activate() | Out-Null
$ws2.Range("G5").Select() | Out-Null

$ws1.Activate() | Out-Null
